# Auto-applied scheduled runner update: refresh cached market-board
# profit figures (currentAveragePrice*, Leve/ItemPrice*, LeveProfit*)
# across the per-job sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 774.25
$ws.Range("I19").Value = 308.9091
$ws.Range("J19").Value = 1168
$ws.Range("K19").Value = 308.9091
$ws.Range("L19").Value = 1168
$ws.Range("M19").Value = -133.9091
$ws.Range("N19").Value = -1518

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H111").Value = 5305.1816
$ws.Range("I111").Value = 3014.5
$ws.Range("J111").Value = 6614.143
$ws.Range("K111").Value = 9043.5
$ws.Range("L111").Value = 19842.429
$ws.Range("M111").Value = -5976.5
$ws.Range("N111").Value = -25976.429

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 2776.5386
$ws.Range("I116").Value = 2293.625
$ws.Range("J116").Value = 3549.2
$ws.Range("K116").Value = 2293.625
$ws.Range("L116").Value = 3549.2
$ws.Range("M116").Value = 1148.375
$ws.Range("N116").Value = -10433.2

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H125").Value = 2550
$ws.Range("I125").Value = 3316.6667
$ws.Range("J125").Value = 250
$ws.Range("K125").Value = 29850.0003
$ws.Range("L125").Value = 2250
$ws.Range("M125").Value = -27390.0003
$ws.Range("N125").Value = -7170

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 2230.25
$ws.Range("I137").Value = 1222.4445
$ws.Range("J137").Value = 2834.9333
$ws.Range("K137").Value = 3667.3335
$ws.Range("L137").Value = 8504.7999
$ws.Range("M137").Value = -1117.3335
$ws.Range("N137").Value = -13604.7999

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 1714.3776
$ws.Range("I138").Value = 685.6667
$ws.Range("J138").Value = 1900.2892
$ws.Range("K138").Value = 2057.0001
$ws.Range("L138").Value = 5700.8676
$ws.Range("M138").Value = 3082.9999
$ws.Range("N138").Value = -15980.8676

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2059.0562
$ws.Range("I32").Value = 2147.3027
$ws.Range("K32").Value = 2147.3027
$ws.Range("M32").Value = -1860.3027

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1777.6154
$ws.Range("I45").Value = 1575.5652
$ws.Range("K45").Value = 1575.5652
$ws.Range("M45").Value = -1198.5652

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 792.38464
$ws.Range("I61").Value = 792.38464
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 792.38464
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -580.38464
$ws.Range("N61").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 944.65515
$ws.Range("I74").Value = 794.2174
$ws.Range("J74").Value = 1521.3334
$ws.Range("K74").Value = 794.2174
$ws.Range("L74").Value = 1521.3334
$ws.Range("M74").Value = 79.7826
$ws.Range("N74").Value = -3269.3334

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 944.65515
$ws.Range("I77").Value = 794.2174
$ws.Range("J77").Value = 1521.3334
$ws.Range("K77").Value = 3971.087
$ws.Range("L77").Value = 7606.666999999999
$ws.Range("M77").Value = 396.913
$ws.Range("N77").Value = -16342.667

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 1792.3539
$ws.Range("I132").Value = 1539.5818
$ws.Range("K132").Value = 4618.7454
$ws.Range("M132").Value = -2088.7454

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 792.38464
$ws.Range("I136").Value = 792.38464
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 2377.15392
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = 172.8460800000003
$ws.Range("N136").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H140").Value = 87613.5
$ws.Range("J140").Value = 87613.5
$ws.Range("L140").Value = 87613.5
$ws.Range("N140").Value = -97973.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 38462890
$ws.Range("I99").Value = 38462890
$ws.Range("K99").Value = 38462890
$ws.Range("M99").Value = -38461392

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1051.8223
$ws.Range("I31").Value = 711.58826
$ws.Range("J31").Value = 2103.4546
$ws.Range("K31").Value = 711.58826
$ws.Range("L31").Value = 2103.4546
$ws.Range("M31").Value = -416.58826
$ws.Range("N31").Value = -2693.4546

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 1051.8223
$ws.Range("I34").Value = 711.58826
$ws.Range("J34").Value = 2103.4546
$ws.Range("K34").Value = 711.58826
$ws.Range("L34").Value = 2103.4546
$ws.Range("M34").Value = -509.58826
$ws.Range("N34").Value = -2507.4546

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 814.85
$ws.Range("I58").Value = 675.25
$ws.Range("J58").Value = 1373.25
$ws.Range("K58").Value = 675.25
$ws.Range("L58").Value = 1373.25
$ws.Range("M58").Value = -472.25
$ws.Range("N58").Value = -1779.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 10242.071
$ws.Range("I132").Value = 12215.909
$ws.Range("J132").Value = 3004.6667
$ws.Range("K132").Value = 36647.727
$ws.Range("L132").Value = 9014.000100000001
$ws.Range("M132").Value = -34117.727
$ws.Range("N132").Value = -14074.0001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 814.85
$ws.Range("I136").Value = 675.25
$ws.Range("J136").Value = 1373.25
$ws.Range("K136").Value = 2025.75
$ws.Range("L136").Value = 4119.75
$ws.Range("M136").Value = 524.25
$ws.Range("N136").Value = -9219.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 13514731
$ws.Range("I131").Value = 111111740
$ws.Range("J131").Value = 1297.7384
$ws.Range("K131").Value = 333335220
$ws.Range("L131").Value = 3893.2152
$ws.Range("M131").Value = -333330180
$ws.Range("N131").Value = -13973.2152

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H133").Value = 3777
$ws.Range("I133").Value = 1757.25
$ws.Range("J133").Value = 4161.7144
$ws.Range("K133").Value = 5271.75
$ws.Range("L133").Value = 12485.1432
$ws.Range("M133").Value = -211.75
$ws.Range("N133").Value = -22605.1432

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H134").Value = 3756.8857
$ws.Range("I134").Value = 1410.6666
$ws.Range("J134").Value = 4981
$ws.Range("K134").Value = 4231.9998
$ws.Range("L134").Value = 14943
$ws.Range("M134").Value = 838.0002000000004
$ws.Range("N134").Value = -25083

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 1700.9429
$ws.Range("I132").Value = 1633.6296
$ws.Range("J132").Value = 1928.125
$ws.Range("K132").Value = 4900.8888
$ws.Range("L132").Value = 5784.375
$ws.Range("M132").Value = -2370.8888
$ws.Range("N132").Value = -10844.375

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1713
$ws.Range("I16").Value = 1581.8334
$ws.Range("J16").Value = 2500
$ws.Range("K16").Value = 1581.8334
$ws.Range("L16").Value = 2500
$ws.Range("M16").Value = -1411.8334
$ws.Range("N16").Value = -2840

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 38627.297
$ws.Range("I132").Value = 1288.9166
$ws.Range("J132").Value = 337334.34
$ws.Range("K132").Value = 3866.7498
$ws.Range("L132").Value = 1012003.02
$ws.Range("M132").Value = -1336.7498
$ws.Range("N132").Value = -1017063.02

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 6444.278
$ws.Range("I136").Value = 6752.7646
$ws.Range("J136").Value = 1200
$ws.Range("K136").Value = 20258.2938
$ws.Range("L136").Value = 3600
$ws.Range("M136").Value = -17708.2938
$ws.Range("N136").Value = -8700

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1594.1698
$ws.Range("I132").Value = 1393.4681
$ws.Range("J132").Value = 3166.3333
$ws.Range("K132").Value = 4180.4043
$ws.Range("L132").Value = 9498.999899999999
$ws.Range("M132").Value = -1650.4043
$ws.Range("N132").Value = -14558.9999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 618.6111
$ws.Range("I136").Value = 324.7619
$ws.Range("K136").Value = 974.2857000000001
$ws.Range("M136").Value = 1575.7143
